# Apply minor adjustment to crash_threshold and related run parameters.

$wb = $excel.ActiveWorkbook

$wsRuns = $wb.Worksheets.Item("runs")
$wsParams = $wb.Worksheets.Item("params")

# --- "runs" sheet updates ---
# Stop Run ID (inclusive): 3 -> 15
$wsRuns.Range("B2").Value = 15
# Sim Time (seconds): 100 -> 500
$wsRuns.Range("B3").Value = 500

# --- "params" sheet updates ---
# Separation Distance (m) column J, rows 2-16: 500 -> 200
for ($r = 2; $r -le 16; $r++) {
    $wsParams.Cells.Item($r, 10).Value = 200
}

# --- Active sheet / selection changes ---
# "params" becomes the active (selected) tab, with selection at N21
$wsParams.Activate()
$wsParams.Range("N21").Select()
